$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two stacked header rows (row 1 + row 2) describing the
# column groups ("Hiver/Eté/Année", "(m3/s)/(MW)/(GWh)"). Replace them
# with a single header row that also introduces two new leading index
# columns (idx / idx2) and explicit date-range / unit labels.
$ws.Rows.Item(1).Delete()

# Row 1 is now the former row 2 (data starts at row 2). Overwrite it with
# the new unified header.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the unit-header cells (F1:K1) the same font used throughout the
# rest of the header/labels.
$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9

$ws.Range("A2:K2").Select()
